$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The lesson-5 wordlist was reorganized: adjectives/na-adjectives, postal
# vocabulary, time words, family/gender words, and food/drink words were
# regrouped into contiguous blocks. Row 1 (header) is unchanged; rows
# 2-127 are rewritten in their new English/Japanese pairing order below.
$ws.Cells.Item(2, 1).Value = 'sea'
$ws.Cells.Item(2, 2).Value = '海|うみ'
$ws.Cells.Item(3, 1).Value = 'postal stamps'
$ws.Cells.Item(3, 2).Value = '切手|きって'
$ws.Cells.Item(4, 1).Value = 'ticket'
$ws.Cells.Item(4, 2).Value = '切符|きっぷ'
$ws.Cells.Item(5, 1).Value = 'surfing'
$ws.Cells.Item(5, 2).Value = 'サーフィン'
$ws.Cells.Item(6, 1).Value = 'homework'
$ws.Cells.Item(6, 2).Value = '宿題|しゅくだい'
$ws.Cells.Item(7, 1).Value = 'food'
$ws.Cells.Item(7, 2).Value = '食べ物|たべもの'
$ws.Cells.Item(8, 1).Value = 'birthday'
$ws.Cells.Item(8, 2).Value = '誕生日|たんじょうび'
$ws.Cells.Item(9, 1).Value = 'test'
$ws.Cells.Item(9, 2).Value = 'テスト'
$ws.Cells.Item(10, 1).Value = 'weather'
$ws.Cells.Item(10, 2).Value = '天気|てんき'
$ws.Cells.Item(11, 1).Value = 'drink'
$ws.Cells.Item(11, 2).Value = '飲み物|のみもの'
$ws.Cells.Item(12, 1).Value = 'postcard'
$ws.Cells.Item(12, 2).Value = '葉書|はがき'
$ws.Cells.Item(13, 1).Value = 'bus'
$ws.Cells.Item(13, 2).Value = 'バス'
$ws.Cells.Item(14, 1).Value = 'airplane'
$ws.Cells.Item(14, 2).Value = '飛行機|ひこうき'
$ws.Cells.Item(15, 1).Value = 'room'
$ws.Cells.Item(15, 2).Value = '部屋|へや'
$ws.Cells.Item(16, 1).Value = 'I (used by men)'
$ws.Cells.Item(16, 2).Value = '僕|ぼく'
$ws.Cells.Item(17, 1).Value = 'holiday; day off; absence'
$ws.Cells.Item(17, 2).Value = '休み|やすみ'
$ws.Cells.Item(18, 1).Value = 'travel'
$ws.Cells.Item(18, 2).Value = '旅行|りょこう'
$ws.Cells.Item(19, 1).Value = 'new'
$ws.Cells.Item(19, 2).Value = '新しい|あたらしい'
$ws.Cells.Item(20, 1).Value = 'hot (weather)'
$ws.Cells.Item(20, 2).Value = '暑い|あつい'
$ws.Cells.Item(21, 1).Value = 'hot (thing)'
$ws.Cells.Item(21, 2).Value = '熱い|あつい'
$ws.Cells.Item(22, 1).Value = 'busy (people/days)'
$ws.Cells.Item(22, 2).Value = '忙しい|いそがしい'
$ws.Cells.Item(23, 1).Value = 'large'
$ws.Cells.Item(23, 2).Value = '大きい|おおきい'
$ws.Cells.Item(24, 1).Value = 'interesting; funny'
$ws.Cells.Item(24, 2).Value = '面白い|おもしろい'
$ws.Cells.Item(25, 1).Value = 'good-looking'
$ws.Cells.Item(25, 2).Value = 'かっこいい'
$ws.Cells.Item(26, 1).Value = 'frightening'
$ws.Cells.Item(26, 2).Value = '怖い|こわい'
$ws.Cells.Item(27, 1).Value = 'cold (weather)'
$ws.Cells.Item(27, 2).Value = '寒い|さむい'
$ws.Cells.Item(28, 1).Value = 'fun'
$ws.Cells.Item(28, 2).Value = '楽しい|たのしい'
$ws.Cells.Item(29, 1).Value = 'small'
$ws.Cells.Item(29, 2).Value = '小さい|ちいさい'
$ws.Cells.Item(30, 1).Value = 'boring'
$ws.Cells.Item(30, 2).Value = 'つまらない'
$ws.Cells.Item(31, 1).Value = 'old (thing)'
$ws.Cells.Item(31, 2).Value = '古い|ふるい'
$ws.Cells.Item(32, 1).Value = 'difficult'
$ws.Cells.Item(32, 2).Value = '難しい|むずかしい'
$ws.Cells.Item(33, 1).Value = 'easy (problem); kind (person)'
$ws.Cells.Item(33, 2).Value = 'やさしい'
$ws.Cells.Item(34, 1).Value = 'inexpensive; cheap (thing)'
$ws.Cells.Item(34, 2).Value = '安い|やすい'
$ws.Cells.Item(35, 1).Value = 'disgusted with; to dislike'
$ws.Cells.Item(35, 2).Value = '嫌い|きらい（な）'
$ws.Cells.Item(36, 1).Value = 'beautiful; clean'
$ws.Cells.Item(36, 2).Value = 'きれい（な）'
$ws.Cells.Item(37, 1).Value = 'healthy; energetic'
$ws.Cells.Item(37, 2).Value = '元気|げんき（な）'
$ws.Cells.Item(38, 1).Value = 'quiet'
$ws.Cells.Item(38, 2).Value = '静か|しずか（な）'
$ws.Cells.Item(39, 1).Value = 'fond of; to like'
$ws.Cells.Item(39, 2).Value = '好き|すき（な）'
$ws.Cells.Item(40, 1).Value = 'to hate'
$ws.Cells.Item(40, 2).Value = '大嫌い|だいきらい（な）'
$ws.Cells.Item(41, 1).Value = 'very fond of; to love'
$ws.Cells.Item(41, 2).Value = '大好き|だいすき（な）'
$ws.Cells.Item(42, 1).Value = 'lively'
$ws.Cells.Item(42, 2).Value = 'にぎやか（な）'
$ws.Cells.Item(43, 1).Value = 'not busy; to have a lot of free time'
$ws.Cells.Item(43, 2).Value = '暇|ひま（な）'
$ws.Cells.Item(44, 1).Value = 'to swim'
$ws.Cells.Item(44, 2).Value = '泳ぐ|およぐ'
$ws.Cells.Item(45, 1).Value = 'to ask'
$ws.Cells.Item(45, 2).Value = '聞く|きく'
$ws.Cells.Item(46, 1).Value = 'to ride; to board'
$ws.Cells.Item(46, 2).Value = '乗る|のる'
$ws.Cells.Item(47, 1).Value = 'to do; to perform'
$ws.Cells.Item(47, 2).Value = 'やる'
$ws.Cells.Item(48, 1).Value = 'to go out'
$ws.Cells.Item(48, 2).Value = '出かける|でかける'
$ws.Cells.Item(49, 1).Value = 'together'
$ws.Cells.Item(49, 2).Value = '一緒に|いっしょに'
$ws.Cells.Item(50, 1).Value = 'extremely'
$ws.Cells.Item(50, 2).Value = 'すごく'
$ws.Cells.Item(51, 1).Value = 'and then'
$ws.Cells.Item(51, 2).Value = 'それから'
$ws.Cells.Item(52, 1).Value = 'It''s okay.; Not to worry.; Everything is under control.'
$ws.Cells.Item(52, 2).Value = '大丈夫|だいじょうぶ'
$ws.Cells.Item(53, 1).Value = 'very'
$ws.Cells.Item(53, 2).Value = 'とても'
$ws.Cells.Item(54, 1).Value = 'what kind of...'
$ws.Cells.Item(54, 2).Value = 'どんな'
$ws.Cells.Item(55, 1).Value = '[counter for flat objects]'
$ws.Cells.Item(55, 2).Value = '～枚|～まい'
$ws.Cells.Item(56, 1).Value = 'to (a place); as far as (a place); till (a time)'
$ws.Cells.Item(56, 2).Value = '～まで'
$ws.Cells.Item(57, 1).Value = 'counter'
$ws.Cells.Item(57, 2).Value = '窓口|まどぐち'
$ws.Cells.Item(58, 1).Value = 'postcard'
$ws.Cells.Item(58, 2).Value = 'はがき'
$ws.Cells.Item(59, 1).Value = 'postal stamps'
$ws.Cells.Item(59, 2).Value = '切手|きって'
$ws.Cells.Item(60, 1).Value = 'parcel'
$ws.Cells.Item(60, 2).Value = '小包|こづつみ'
$ws.Cells.Item(61, 1).Value = 'airmail'
$ws.Cells.Item(61, 2).Value = '航空便|こうくうびん'
$ws.Cells.Item(62, 1).Value = 'surface mail'
$ws.Cells.Item(62, 2).Value = '船便|ふなびん'
$ws.Cells.Item(63, 1).Value = 'special delivery'
$ws.Cells.Item(63, 2).Value = '速達|そくたつ'
$ws.Cells.Item(64, 1).Value = 'registered mail'
$ws.Cells.Item(64, 2).Value = '書留|かきとめ'
$ws.Cells.Item(65, 1).Value = 'insurance'
$ws.Cells.Item(65, 2).Value = '保険|ほけん'
$ws.Cells.Item(66, 1).Value = 'Can you take care of this, please?'
$ws.Cells.Item(66, 2).Value = 'これ、お願いします。'
$ws.Cells.Item(67, 1).Value = 'Give me three 50-yen stamps, please.'
$ws.Cells.Item(67, 2).Value = '五十円切手を三枚ください。'
$ws.Cells.Item(68, 1).Value = 'Make this (an airmail), please.'
$ws.Cells.Item(68, 2).Value = '(航空便)でお願いします。'
$ws.Cells.Item(69, 1).Value = 'How many days will it take?'
$ws.Cells.Item(69, 2).Value = '何日ぐらいかかりますか。'
$ws.Cells.Item(70, 1).Value = 'It will be 150 yen.'
$ws.Cells.Item(70, 2).Value = '百五十円になります。'
$ws.Cells.Item(71, 1).Value = 'mountain'
$ws.Cells.Item(71, 2).Value = '山|やま'
$ws.Cells.Item(72, 1).Value = 'Mr./Ms. Yamakawa'
$ws.Cells.Item(72, 2).Value = '山川さん|やまかわさん'
$ws.Cells.Item(73, 1).Value = 'Mt. Fuji'
$ws.Cells.Item(73, 2).Value = '富士山|ふじさん'
$ws.Cells.Item(74, 1).Value = 'river'
$ws.Cells.Item(74, 2).Value = '川|かわ'
$ws.Cells.Item(75, 1).Value = 'Mr./Ms. Ogawa'
$ws.Cells.Item(75, 2).Value = '小川さん|おがわさん'
$ws.Cells.Item(76, 1).Value = 'fine'
$ws.Cells.Item(76, 2).Value = '元気な|げんきな'
$ws.Cells.Item(77, 1).Value = 'the first day of the year'
$ws.Cells.Item(77, 2).Value = '元日|がんじつ'
$ws.Cells.Item(78, 1).Value = 'local'
$ws.Cells.Item(78, 2).Value = '地元|じもと'
$ws.Cells.Item(79, 1).Value = 'fine'
$ws.Cells.Item(79, 2).Value = '元気な|げんきな'
$ws.Cells.Item(80, 1).Value = 'weather'
$ws.Cells.Item(80, 2).Value = '天気|てんき'
$ws.Cells.Item(81, 1).Value = 'electricity'
$ws.Cells.Item(81, 2).Value = '電気|でんき'
$ws.Cells.Item(82, 1).Value = 'feeling'
$ws.Cells.Item(82, 2).Value = '気持ち|きもち'
$ws.Cells.Item(83, 1).Value = 'popularity'
$ws.Cells.Item(83, 2).Value = '人気|にんき'
$ws.Cells.Item(84, 1).Value = 'heaven'
$ws.Cells.Item(84, 2).Value = '天国|てんごく'
$ws.Cells.Item(85, 1).Value = 'Japanese emperor'
$ws.Cells.Item(85, 2).Value = '天皇|てんのう'
$ws.Cells.Item(86, 1).Value = 'genius'
$ws.Cells.Item(86, 2).Value = '天才|てんさい'
$ws.Cells.Item(87, 1).Value = 'I'
$ws.Cells.Item(87, 2).Value = '私|わたし'
$ws.Cells.Item(88, 1).Value = 'private university'
$ws.Cells.Item(88, 2).Value = '私立大学|しりつだいがく'
$ws.Cells.Item(89, 1).Value = 'private railroad'
$ws.Cells.Item(89, 2).Value = '私鉄|してつ'
$ws.Cells.Item(90, 1).Value = 'now'
$ws.Cells.Item(90, 2).Value = '今|いま'
$ws.Cells.Item(91, 1).Value = 'today'
$ws.Cells.Item(91, 2).Value = '今日|きょう'
$ws.Cells.Item(92, 1).Value = 'tonight'
$ws.Cells.Item(92, 2).Value = '今晩|こんばん'
$ws.Cells.Item(93, 1).Value = 'this month'
$ws.Cells.Item(93, 2).Value = '今月|こんげつ'
$ws.Cells.Item(94, 1).Value = 'this year'
$ws.Cells.Item(94, 2).Value = '今年|ことし'
$ws.Cells.Item(95, 1).Value = 'Mr./Ms. Tanaka'
$ws.Cells.Item(95, 2).Value = '田中さん|たなかさん'
$ws.Cells.Item(96, 1).Value = 'Mr./Ms. Yamada'
$ws.Cells.Item(96, 2).Value = '山田さん|やまださん'
$ws.Cells.Item(97, 1).Value = 'rice field'
$ws.Cells.Item(97, 2).Value = '田んぼ|たんぼ'
$ws.Cells.Item(98, 1).Value = 'woman (kun-yomi)'
$ws.Cells.Item(98, 2).Value = '女の人|おんなのひと'
$ws.Cells.Item(99, 1).Value = 'woman (on-yomi)'
$ws.Cells.Item(99, 2).Value = '女性|じょせい'
$ws.Cells.Item(100, 1).Value = 'girl'
$ws.Cells.Item(100, 2).Value = '女の子|おんなのこ'
$ws.Cells.Item(101, 1).Value = 'the eldest daughter'
$ws.Cells.Item(101, 2).Value = '長女|ちょうじょ'
$ws.Cells.Item(102, 1).Value = 'man (kun-yomi)'
$ws.Cells.Item(102, 2).Value = '男の人|おとこのひと'
$ws.Cells.Item(103, 1).Value = 'man (on-yomi)'
$ws.Cells.Item(103, 2).Value = '男性|だんせい'
$ws.Cells.Item(104, 1).Value = 'boy'
$ws.Cells.Item(104, 2).Value = '男の子|おとこのこ'
$ws.Cells.Item(105, 1).Value = 'male student'
$ws.Cells.Item(105, 2).Value = '男子学生|だんしがくせい'
$ws.Cells.Item(106, 1).Value = 'to see'
$ws.Cells.Item(106, 2).Value = '見る|みる'
$ws.Cells.Item(107, 1).Value = 'sightseeing'
$ws.Cells.Item(107, 2).Value = '見物|けんぶつ'
$ws.Cells.Item(108, 1).Value = 'flower viewing'
$ws.Cells.Item(108, 2).Value = '花見|はなみ'
$ws.Cells.Item(109, 1).Value = 'opinion'
$ws.Cells.Item(109, 2).Value = '意見|いけん'
$ws.Cells.Item(110, 1).Value = 'to go'
$ws.Cells.Item(110, 2).Value = '行く|いく'
$ws.Cells.Item(111, 1).Value = 'bank'
$ws.Cells.Item(111, 2).Value = '銀行|ぎんこう'
$ws.Cells.Item(112, 1).Value = 'first line'
$ws.Cells.Item(112, 2).Value = '一行目|いちぎょうめ'
$ws.Cells.Item(113, 1).Value = 'travel'
$ws.Cells.Item(113, 2).Value = '旅行|りょこう'
$ws.Cells.Item(114, 1).Value = 'to eat'
$ws.Cells.Item(114, 2).Value = '食べる|たべる'
$ws.Cells.Item(115, 1).Value = 'food'
$ws.Cells.Item(115, 2).Value = '食べ物|たべもの'
$ws.Cells.Item(116, 1).Value = 'cafeteria'
$ws.Cells.Item(116, 2).Value = '食堂|しょくどう'
$ws.Cells.Item(117, 1).Value = 'meal'
$ws.Cells.Item(117, 2).Value = '食事|しょくじ'
$ws.Cells.Item(118, 1).Value = 'breakfast'
$ws.Cells.Item(118, 2).Value = '朝食|ちょうしょく'
$ws.Cells.Item(119, 1).Value = 'to drink'
$ws.Cells.Item(119, 2).Value = '飲む|のむ'
$ws.Cells.Item(120, 1).Value = 'drink'
$ws.Cells.Item(120, 2).Value = '飲み物|のみもの'
$ws.Cells.Item(121, 1).Value = 'drunken driving'
$ws.Cells.Item(121, 2).Value = '飲酒運転|いんしゅうんてん'
$ws.Cells.Item(122, 1).Value = 'coffee'
$ws.Cells.Item(122, 2).Value = 'コーヒー'
$ws.Cells.Item(123, 1).Value = 'concert'
$ws.Cells.Item(123, 2).Value = 'コンサート'
$ws.Cells.Item(124, 1).Value = 'Vienna'
$ws.Cells.Item(124, 2).Value = 'ウィーン'
$ws.Cells.Item(125, 1).Value = 'cafe'
$ws.Cells.Item(125, 2).Value = 'カフェ'
$ws.Cells.Item(126, 1).Value = 'classical music'
$ws.Cells.Item(126, 2).Value = 'クラシック'
$ws.Cells.Item(127, 1).Value = 'cake'
$ws.Cells.Item(127, 2).Value = 'ケーキ'
